# New MR82F001 Coremark results: six new rows (48-53) appended after the
# existing STM32F401 block, plus the supporting formatting tweaks that Excel
# applied when the data was typed in (quote-prefixed optimisation labels in
# column C, a wider column C, and the selection/scroll position left on the
# last new cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 48: MCU label + first data row ------------------------------------
$ws.Cells.Item(48, 1).Value = "MR82F001"
$ws.Cells.Item(48, 2).Value = 168
$ws.Cells.Item(48, 3).Value = "'AC6 -Ofast -Otime -Omax Icache + Dcache"
$ws.Cells.Item(48, 4).Value = 544.61
$ws.Cells.Item(48, 5).Formula = "=D48/B48"

# --- Row 49 -----------------------------------------------------------------
$ws.Cells.Item(49, 2).Value = 168
$ws.Cells.Item(49, 3).Value = "'AC6 -Ofast -Otime Icache+Dcache"
$ws.Cells.Item(49, 4).Value = 480.92
$ws.Cells.Item(49, 5).Formula = "=D49/B49"

# --- Row 50 -----------------------------------------------------------------
$ws.Cells.Item(50, 2).Value = 168
$ws.Cells.Item(50, 3).Value = "'AC6 -Ofast Icache+Dcache"
$ws.Cells.Item(50, 4).Value = 492.44
$ws.Cells.Item(50, 5).Formula = "=D50/B50"

# --- Row 51 -----------------------------------------------------------------
$ws.Cells.Item(51, 2).Value = 168
$ws.Cells.Item(51, 3).Value = "'AC6 -Ofast  -Otime -Omax Icache"
$ws.Cells.Item(51, 4).Value = 397.61

# --- Row 52 -----------------------------------------------------------------
$ws.Cells.Item(52, 2).Value = 168
$ws.Cells.Item(52, 3).Value = "'AC6 -Ofast -Otime Icache"
$ws.Cells.Item(52, 4).Value = 338.94

# --- Row 53 -----------------------------------------------------------------
$ws.Cells.Item(53, 2).Value = 168
$ws.Cells.Item(53, 3).Value = "'AC6 -Ofast Icache"
$ws.Cells.Item(53, 4).Value = 342.62

# Rows 51-53 got their E column filled as one drag-fill, i.e. one shared
# formula spanning the three cells (mirrors si="0"/"1"/"2" already in the
# sheet).
$ws.Range("E51:E53").Formula = "=D51/B51"

# Column C needed to grow to fit the new, longer optimisation strings.
$ws.Columns.Item(3).ColumnWidth = 36.45

# Leave the selection/scroll where the user ended up after typing the data.
$ws.Cells.Item(53, 4).Select()
$excel.ActiveWindow.ScrollRow = 29
$excel.ActiveWindow.ScrollColumn = 1
